$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (extr1..extr8) with new C/D/E values ---

# Row 8 (extr1): C 5->14, D 12->11, E 0->1
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = 11
$ws.Cells.Item(8,5).Value = $true

# Row 9 (extr2): C 5->16, E 0->1
$ws.Cells.Item(9,3).Value = 16
$ws.Cells.Item(9,5).Value = $true

# Row 10 (extr3): C 10->5, D 11->12, E 0->1
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 12
$ws.Cells.Item(10,5).Value = $true

# Row 11 (extr4): C 7->5, D 8->9
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 9

# Row 12 (extr5): C 9->10
$ws.Cells.Item(12,3).Value = 10

# Row 13 (extr6): D 11->8
$ws.Cells.Item(13,4).Value = 8

# Row 14 (extr7): C 5->9, D 7->11, E 1->0
$ws.Cells.Item(14,3).Value = 9
$ws.Cells.Item(14,4).Value = 11
$ws.Cells.Item(14,5).Value = $false

# Row 15 (extr8): C 8->7, D 5->11, E 1->0
$ws.Cells.Item(15,3).Value = 7
$ws.Cells.Item(15,4).Value = 11
$ws.Cells.Item(15,5).Value = $false

# --- Append new rows 16 (line7) and 17 (line8), copying formatting from row 15's A column ---
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "line7"
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 7
$ws.Cells.Item(16,5).Value = $true

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "line8"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = $false
